# define fill up event , marker_pen test
#
# Adds two new "test case" blocks (selector / pen) to the STB_marker sheet,
# mirroring the layout already used on the STB_timer sheet (title cell in
# column E with the highlighted "s=4" fill, then a numbered F/G step list
# terminated by an "END" marker in column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STB_marker")
$wsTimer = $wb.Worksheets.Item("STB_timer")

# Style donors already present in the workbook:
#  - STB_timer!E2 carries the highlighted title style (s="4")
#  - STB_marker!C1 carries the plain body style (s="2")
$titleStyleSrc = $wsTimer.Range("E2")
$bodyStyleSrc = $ws.Range("C1")

function Set-StepCell {
    param($range, $value, $styleSrc)
    $range.Value = $value
    $styleSrc.Copy()
    $range.PasteSpecial(-4122)
}

function Touch-RowSpan {
    param($row)
    # Lightly "touch" the outer A/K columns (same as the template rows) so
    # the row's computed column span matches the rest of the table, without
    # leaving any real cell content behind.
    $ws.Range("A$row").Value = ""
    $ws.Range("K$row").Value = ""
}

# ---------------------------------------------------------------------
# Block 1: "selector" test case -> rows 2-12
# ---------------------------------------------------------------------
$selectorSteps = @(
    "open STB",
    "STB all tools button",
    "open Marker",
    "move bar right",
    "move bar up",
    "select selector mode",
    "tap 'myViewBoard Display' on hotseat",
    "current app compare",
    "close button",
    "homepage"
)

Touch-RowSpan 2
Set-StepCell $ws.Range("E2") "selector" $titleStyleSrc

for ($i = 0; $i -lt $selectorSteps.Count; $i++) {
    $row = 2 + $i
    $stepNum = $i + 1
    if ($row -ne 2) { Touch-RowSpan $row }
    Set-StepCell $ws.Range("F$row") $stepNum $bodyStyleSrc
    Set-StepCell $ws.Range("G$row") $selectorSteps[$i] $bodyStyleSrc
}

Touch-RowSpan 12
Set-StepCell $ws.Range("G12") "END" $bodyStyleSrc

# ---------------------------------------------------------------------
# Block 2: "pen" test case -> rows 14-25
# ---------------------------------------------------------------------
$penSteps = @(
    "open STB",
    "STB all tools button",
    "open Marker",
    "move bar right",
    "move bar up",
    "select pen mode",
    "Screenshot compare_1",
    "fill up left upper coner by 50 steps",
    "Screenshot compare_2",
    "compare different",
    "close button"
)

Touch-RowSpan 14
Set-StepCell $ws.Range("E14") "pen" $titleStyleSrc

for ($i = 0; $i -lt $penSteps.Count; $i++) {
    $row = 14 + $i
    $stepNum = $i + 1
    if ($row -ne 14) { Touch-RowSpan $row }
    Set-StepCell $ws.Range("F$row") $stepNum $bodyStyleSrc
    Set-StepCell $ws.Range("G$row") $penSteps[$i] $bodyStyleSrc
}

Set-StepCell $ws.Range("G25") "END" $bodyStyleSrc

# ---------------------------------------------------------------------
# View/selection tweaks
# ---------------------------------------------------------------------

# STB_timer keeps its own remembered selection even though it's not the
# active tab.
$wsTimer.Activate()
$wsTimer.Range("K6").Select()

# STB_marker stays the active tab, scrolled down a bit with a new
# selection, matching the edited state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H20").Select()
